# Edit script: append new scrape batch (2025-12-24 18:27:26 JST) to the
# "案件情報" workbook's ランサーズ (Lancers) sheet.
#
# Behaviour mirrors the upstream scraper: the freshest 3 postings are
# merged to the top/middle of the existing list (by priority score) and
# every row's "取得日時" timestamp is refreshed to the latest run time.
# The net result is 3 additional rows (12 -> 15) with updated content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$data = @(
    @{A='2025-12-24 18:27:26'; B='製造業向け図面自動生成システムの開発・ツール化を支援してくださるエンジニア募集(AI/バックエンド)'; C='システム開発'; D='300,000 円 ~ 500,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5460562'; G=435; H='🔥AI,Ai ◆ツール,開発'}
    @{A='2025-12-24 18:27:26'; B='産業機械向けAI異常検知・状態推定システムの開発・導入支援エンジニア募集(AI/エッジ・組み込み)'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5450864'; G=383; H='🔥AI,Ai ◆開発'}
    @{A='2025-12-24 18:27:26'; B='【急募】自社AIプロダクト開発|バックエンドエンジニア'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5460544'; G=375; H='🔥AI,Ai ◆開発'}
    @{A='2025-12-24 18:27:26'; B='【フルタイム】最先端AI(LLM)開発エンジニア募集!新規プロダクトの核となる開発パートナーを募集'; C='システム開発'; D='1,000,000 円 ~ 3,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5460294'; G=375; H='🔥AI,Ai ◆開発'}
    @{A='2025-12-24 18:27:26'; B='【週5日】法人向け生成AIサービス(RAG・議事録機能)のコア開発を担うリードエンジニア募集'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5460267'; G=375; H='🔥AI,Ai ◆開発'}
    @{A='2025-12-24 18:27:26'; B='【急募】AI活用でPDFタイトル修正のフリーランス募集!'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5459721'; G=310; H='🔥AI,Ai'}
    @{A='2025-12-24 18:27:26'; B='施設管理・現場業務向け チェックリスト業務の自動化・報告書作成システム開発エンジニア募集'; C='システム開発'; D='300,000 円 ~ 500,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5460563'; G=220; H='◆開発,システム開発 ◇管理'}
    @{A='2025-12-24 18:27:26'; B='【急募】野球スコアボードシステム開発のフリーランス募集'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5459984'; G=118; H='◆開発,システム開発'}
    @{A='2025-12-24 18:27:26'; B='初回 【AWSクラウドリフト】業務アプリ移行支援エンジニア募集(Java / .NET)'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5459847'; G=103; H='★Java ◇アプリ'}
    @{A='2025-12-24 18:27:26'; B='現品票管理・納品書・請求書のシステムづくり'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5459942'; G=53; H='◇管理'}
    @{A='2025-12-24 18:27:26'; B='急募 限定公開 限定公開の仕事'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5460299'; G=18; H=$null}
    @{A='2025-12-24 18:27:26'; B='【電卓設計】ハードウェアとソフトウェアの専門家を募集!'; C='システム開発'; D='50,000 円 ~ 100,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5459773'; G=18; H=$null}
    @{A='2025-12-24 18:27:26'; B='【電卓設計】ハードウェアとソフトウェアの専門家を募集!'; C='システム開発'; D='50,000 円 ~ 100,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5459232'; G=18; H=$null}
    @{A='2025-12-24 18:27:26'; B='【急募】お名前VPSでのFTP・WPファイルアップロード改善依頼'; C='システム開発'; D='5,000 円 ~ 10,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5459964'; G=10; H=$null}
)

# Remove every existing hyperlink on the sheet first. (In this COM
# implementation Range.Hyperlinks.Delete() clears the whole worksheet's
# hyperlink collection, so one call up-front is sufficient and avoids
# leaving stale/duplicate relationships behind when URLs change.)
$ws.Range("A1").Hyperlinks.Delete()

# Write out the full 14-row data block (rows 2-15) and re-create the
# hyperlink on each URL cell in column F.
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $item = $data[$i]

    $ws.Cells.Item($row, 1).Value = $item.A
    $ws.Cells.Item($row, 2).Value = $item.B
    $ws.Cells.Item($row, 3).Value = $item.C
    $ws.Cells.Item($row, 4).Value = $item.D
    $ws.Cells.Item($row, 5).Value = $item.E

    $fcell = $ws.Cells.Item($row, 6)
    $fcell.Value = $item.F
    $ws.Hyperlinks.Add($fcell, $item.F) | Out-Null

    $ws.Cells.Item($row, 7).Value = $item.G

    if ($item.H -ne $null) {
        $ws.Cells.Item($row, 8).Value = $item.H
    }
}

# Column width tweaks (B: 51 -> 52, H: 12 -> 16 characters). The COM
# ColumnWidth property is offset from the stored OOXML column width by
# the sheet's default ~0.83 character padding, so subtract that back out
# to land on the exact target widths.
$ws.Columns.Item(2).ColumnWidth = 51.17
$ws.Columns.Item(8).ColumnWidth = 15.17
